# Updates the "cryptos" list (Sheet1) with refreshed price / 1h-volume data,
# matching the GitHub Actions scheduled refresh commit.
#
# All data cells in columns B:E are stored as *text* (t="inlineStr" in the
# source workbook), even though many of the Price values look numeric
# (e.g. "1.00", "0.410", "95.527.43"). A plain `Range.Value = "..."` write
# lets Excel's automatic type inference kick in and silently turn those
# number-looking strings into real numbers (losing trailing zeros / the
# thousand-dot formatting, e.g. "1.00" -> 1, "0.410" -> 0.41), and it also
# happens to stamp a new NumberFormat style onto the cell. To avoid both
# problems we:
#   1. write the literal value as a quoted-string formula ( ="text" ),
#      which forces a text result regardless of what the text looks like;
#   2. Copy the cell and PasteSpecial *values only* back onto itself,
#      collapsing the formula down to a plain literal text value again
#      (no "<f>" left behind, no style change).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($rangeAddr, $text) {
    $escapedQuotes = $text -replace '"', '""'
    $ws.Range($rangeAddr).Formula = '="' + $escapedQuotes + '"'
    $ws.Range($rangeAddr).Copy()
    $ws.Range($rangeAddr).PasteSpecial(-4163)  # xlPasteValues
}

$excel.CutCopyMode = 0

# --- Row 2: Bitcoin ---
Set-TextValue 'D2' '95.527.43'
Set-TextValue 'E2' '  +3.11%  '

# --- Row 3: Ethereum ---
Set-TextValue 'D3' '3.603.96'
Set-TextValue 'E3' '  +7.15%  '

# --- Row 4: TetherUSD ---
Set-TextValue 'E4' '  -0.04%  '

# --- Row 5: Solana ---
Set-TextValue 'D5' '240.42'
Set-TextValue 'E5' '  +3.73%  '

# --- Row 6: BNB ---
Set-TextValue 'D6' '652.09'
Set-TextValue 'E6' '  +5.68%  '

# --- Row 7: XRP ---
Set-TextValue 'E7' '  +7.09%  '

# --- Row 8: Dogecoin ---
Set-TextValue 'D8' '0.410'
Set-TextValue 'E8' '  +5.79%  '

# --- Row 9: USDC ---
Set-TextValue 'E9' '  -0.13%  '

# --- Row 10: Cardano ---
Set-TextValue 'D10' '1.00'
Set-TextValue 'E10' '  +5.84%  '

# --- Row 11: LidoStakedEther ---
Set-TextValue 'D11' '3.601.04'
Set-TextValue 'E11' '  +7.13%  '

# --- Row 12: Avalanche ---
Set-TextValue 'D12' '43.21'
Set-TextValue 'E12' '  +1.61%  '

# --- Row 13: TRON ---
Set-TextValue 'D13' '0.200'
Set-TextValue 'E13' '  +1.86%  '

# --- Row 14: Toncoin ---
Set-TextValue 'D14' '6.32'
Set-TextValue 'E14' '  +2.10%  '

# --- Row 15: WrappedliquidstakedEther2.0 ---
Set-TextValue 'D15' '4.284.44'
Set-TextValue 'E15' '  +7.32%  '

# --- Row 16: WrappedBTC ---
Set-TextValue 'D16' '95.370.35'
Set-TextValue 'E16' '  +3.04%  '

# --- Row 17: ShibaInu ---
Set-TextValue 'D17' '0.0000257'
Set-TextValue 'E17' '  +5.46%  '

# --- Row 18: WrappedEther ---
Set-TextValue 'D18' '3.599.10'
Set-TextValue 'E18' '  +7.21%  '

# --- Row 19: Polkadot ---
Set-TextValue 'D19' '7.95'
Set-TextValue 'E19' '  -1.36%  '

# --- Row 20: Uniswap ---
Set-TextValue 'D20' '12.51'
Set-TextValue 'E20' '  +11.03%  '

# --- Row 21: Chainlink ---
Set-TextValue 'D21' '18.06'
Set-TextValue 'E21' '  +3.80%  '

# --- Row 22: SuiNetwork ---
Set-TextValue 'D22' '3.49'
Set-TextValue 'E22' '  +5.02%  '

# --- Row 23: Stellar ---
Set-TextValue 'E23' '  +13.03%  '

# --- Row 24: BitcoinCash ---
Set-TextValue 'D24' '512.13'
Set-TextValue 'E24' '  +3.73%  '

# --- Row 25: PEPE ---
Set-TextValue 'D25' '0.0000196'
Set-TextValue 'E25' '  +8.14%  '

# --- Row 26: NEARProtocol ---
Set-TextValue 'D26' '6.68'
Set-TextValue 'E26' '  +2.12%  '

# --- Row 27: Litecoin ---
Set-TextValue 'D27' '96.38'
Set-TextValue 'E27' '  +6.63%  '

# --- Row 28: Aptos ---
Set-TextValue 'D28' '12.93'
Set-TextValue 'E28' '  +8.31%  '

# --- Row 29: PancakeSwap ---
Set-TextValue 'D29' '3.13'
Set-TextValue 'E29' '  +17.15%  '

# --- Row 30: InternetComputer(DFINITY) ---
Set-TextValue 'D30' '11.31'
Set-TextValue 'E30' '  +2.56%  '

# --- Row 31: Dai ---
Set-TextValue 'D31' '1.00'
Set-TextValue 'E31' '  -0.06%  '

# --- Row 32: Hedera ---
Set-TextValue 'E32' '  +2.68%  '

# --- Row 33: Binance-PegBSC-USD ---
Set-TextValue 'D33' '0.992'
Set-TextValue 'E33' '  -0.97%  '

# --- Row 34: Cronos ---
Set-TextValue 'D34' '0.176'
Set-TextValue 'E34' '  +2.95%  '

# --- Row 35: EthereumClassic ---
Set-TextValue 'D35' '31.92'
Set-TextValue 'E35' '  +11.95%  '

# --- Row 36: PolygonEcosystemToken ---
Set-TextValue 'D36' '0.560'
Set-TextValue 'E36' '  +6.61%  '

# --- Row 37: RenderToken ---
Set-TextValue 'D37' '8.28'
Set-TextValue 'E37' '  +11.50%  '

# --- Row 38: Bittensor ---
Set-TextValue 'D38' '561.73'
Set-TextValue 'E38' '  +2.03%  '

# --- Row 39: Fetch.AI ---
Set-TextValue 'D39' '1.47'
Set-TextValue 'E39' '  +7.18%  '

# --- Row 40: USDe ---
Set-TextValue 'E40' '  -0.09%  '

# --- Row 41: ARBITRUM ---
Set-TextValue 'D41' '0.930'
Set-TextValue 'E41' '  +5.74%  '

# --- Row 42: Kaspa ---
Set-TextValue 'D42' '0.151'
Set-TextValue 'E42' '  +1.22%  '

# --- Row 43: ImmutableX ---
Set-TextValue 'E43' '  +2.75%  '

# --- Row 44/45: WhiteBITCoin <-> Filecoin swap places (rank reordering) ---
Set-TextValue 'B44' 'Filecoin'
Set-TextValue 'C44' 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextValue 'D44' '5.71'
Set-TextValue 'E44' '  +6.29%  '

Set-TextValue 'B45' 'WhiteBITCoin'
Set-TextValue 'C45' 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
Set-TextValue 'D45' '23.76'
Set-TextValue 'E45' '  +0.43%  '

# --- Row 46: Stacks ---
Set-TextValue 'D46' '2.27'
Set-TextValue 'E46' '  +8.82%  '

# --- Row 47: VeChain ---
Set-TextValue 'D47' '0.0419'
Set-TextValue 'E47' '  +3.98%  '

# --- Row 48/49: OKB <-> EnergySwap swap places (rank reordering) ---
Set-TextValue 'B48' 'EnergySwap'
Set-TextValue 'C48' 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue 'D48' '33.26'
Set-TextValue 'E48' '  +42.84%  '

Set-TextValue 'B49' 'OKB'
Set-TextValue 'C49' 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
Set-TextValue 'D49' '54.35'
Set-TextValue 'E49' '  +3.43%  '

# --- Row 50: MantraDAO ---
Set-TextValue 'D50' '3.45'
Set-TextValue 'E50' '  -3.74%  '

# --- Row 51: Cosmos ---
Set-TextValue 'D51' '8.10'
Set-TextValue 'E51' '  +2.98%  '

$excel.CutCopyMode = 0
